$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell-value updates from the "Updated cryptos list" data refresh.
# D-column (Price) cells often look numeric ("14.10", "1.00", ...). A plain
# Range.Value assignment would let Excel auto-coerce them to numbers and
# silently drop significant trailing zeros / thousands separators, so for
# every D-column write we briefly force a Text number format, assign the
# literal string, then ClearFormats() so the cells style reverts to the
# workbook default (matching the original, unstyled cells).

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "39.348.29"
$c.ClearFormats()
$ws.Range("E2").Value = "  -3.21%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.216.31"
$c.ClearFormats()
$ws.Range("E3").Value = "  -6.59%  "
$ws.Range("E4").Value = "  +0.05%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "296.35"
$c.ClearFormats()
$ws.Range("E5").Value = "  -4.71%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "83.07"
$c.ClearFormats()
$ws.Range("E6").Value = "  -4.21%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.510"
$c.ClearFormats()
$ws.Range("E7").Value = "  -3.71%  "
$ws.Range("E8").Value = "  +0.03%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.467"
$c.ClearFormats()
$ws.Range("E9").Value = "  -4.80%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.0774"
$c.ClearFormats()
$ws.Range("E10").Value = "  -7.92%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "29.08"
$c.ClearFormats()
$ws.Range("E11").Value = "  -4.35%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "47.85"
$c.ClearFormats()
$ws.Range("E12").Value = "  -9.28%  "
$ws.Range("E13").Value = "  -2.16%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "2.567.18"
$c.ClearFormats()
$ws.Range("E14").Value = "  -5.90%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "6.30"
$c.ClearFormats()
$ws.Range("E15").Value = "  -3.45%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "14.10"
$c.ClearFormats()
$ws.Range("E16").Value = "  -5.77%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "2.216.32"
$c.ClearFormats()
$ws.Range("E17").Value = "  -7.38%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.714"
$c.ClearFormats()
$ws.Range("E18").Value = "  -5.71%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "39.260.39"
$c.ClearFormats()
$ws.Range("E19").Value = "  -3.25%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.0₃0873"
$c.ClearFormats()
$ws.Range("E20").Value = "  -4.04%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "5.70"
$c.ClearFormats()
$ws.Range("E21").Value = "  -6.91%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "64.95"
$c.ClearFormats()
$ws.Range("E22").Value = "  -5.26%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "10.26"
$c.ClearFormats()
$ws.Range("E23").Value = "  -4.62%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "227.76"
$c.ClearFormats()
$ws.Range("E24").Value = "  -3.23%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.ClearFormats()
$ws.Range("E25").Value = "  -0.22%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "2.40"
$c.ClearFormats()
$ws.Range("E26").Value = "  -6.72%  "
$ws.Range("E27").Value = "  -0.08%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "22.51"
$c.ClearFormats()
$ws.Range("E28").Value = "  -5.42%  "
$ws.Range("E29").Value = "  -2.85%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "9.10"
$c.ClearFormats()
$ws.Range("E30").Value = "  -1.55%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "149.81"
$c.ClearFormats()
$ws.Range("E31").Value = "  -2.67%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "31.91"
$c.ClearFormats()
$ws.Range("E32").Value = "  -6.72%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.ClearFormats()
$ws.Range("E33").Value = "  -0.11%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "4.84"
$c.ClearFormats()
$ws.Range("E34").Value = "  -6.61%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.0693"
$c.ClearFormats()
$ws.Range("E35").Value = "  -4.84%  "
$ws.Range("E36").Value = "  -3.37%  "
$ws.Range("E37").Value = "  -3.85%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.0962"
$c.ClearFormats()
$ws.Range("E38").Value = "  -3.99%  "
$ws.Range("B39").Value = "Celestia"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "15.14"
$c.ClearFormats()
$ws.Range("E39").Value = "  -5.34%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "2.63"
$c.ClearFormats()
$ws.Range("E40").Value = "  -4.69%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "1.64"
$c.ClearFormats()
$ws.Range("E41").Value = "  -3.76%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "3.64"
$c.ClearFormats()
$ws.Range("E42").Value = "  -5.35%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "1.907.76"
$c.ClearFormats()
$ws.Range("E43").Value = "  -2.56%  "
$ws.Range("E44").Value = "  -3.72%  "
$ws.Range("E45").Value = "  -15.80%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "16.03"
$c.ClearFormats()
$ws.Range("E46").Value = "  -8.50%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "8.93"
$c.ClearFormats()
$ws.Range("E47").Value = "  -4.70%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "2.61"
$c.ClearFormats()
$ws.Range("E48").Value = "  -2.65%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "2.432.26"
$c.ClearFormats()
$ws.Range("E49").Value = "  -6.33%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "70.28"
$c.ClearFormats()
$ws.Range("E50").Value = "  -2.16%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "86.77"
$c.ClearFormats()
$ws.Range("E51").Value = "  -6.74%  "
